# Variable_Name_Metadata.xlsx — drop the phosphorus/nitrogen driver-variable
# rows (OGM_don, NIT_nit, NIT_amm, PHS_frp) from the "inflow.csv" block now
# that those drivers are produced by the new R script instead of by hand in
# Excel. That block was rows 14-17 (the PHS_frp row plus the three others
# sharing its formatting); the trailing "outflow.csv" example (old rows
# 18-19) slides up to take rows 14-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four now-obsolete variable rows; everything below shifts up,
# which is also what renumbers/prunes the shared-string table and the used
# range on save.
$ws.Rows("14:17").Delete() | Out-Null

# Leave the selection where the author's last click landed.
$ws.Range("E13").Select() | Out-Null
